$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDA")
$ws.Range("I19:J29").Select()
$formula = "=LET(d,TAKE(SORT(C3:G12,{2,3,4},{-1,-1,-1}),,1),VSTACK(L2:M2,HSTACK(SEQUENCE(ROWS(d)),d)))"
$excel.Selection.FormulaArray = $formula
Write-Output "set"
